$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 38785.332  # H21
$ws.Cells.Item(21, 10).Value = 32538.6  # J21
$ws.Cells.Item(21, 12).Value = 32538.6  # L21
$ws.Cells.Item(21, 14).Value = -33474.6  # N21
$ws.Cells.Item(23, 8).Value = 38785.332  # H23
$ws.Cells.Item(23, 10).Value = 32538.6  # J23
$ws.Cells.Item(23, 12).Value = 32538.6  # L23
$ws.Cells.Item(23, 14).Value = -33006.6  # N23
$ws.Cells.Item(26, 8).Value = 44994.5  # H26
$ws.Cells.Item(26, 10).Value = 44994.5  # J26
$ws.Cells.Item(26, 12).Value = 44994.5  # L26
$ws.Cells.Item(26, 14).Value = -45682.5  # N26
$ws.Cells.Item(39, 8).Value = 188.2  # H39
$ws.Cells.Item(39, 9).Value = 85.25  # I39
$ws.Cells.Item(39, 10).Value = 600  # J39
$ws.Cells.Item(39, 11).Value = 255.75  # K39
$ws.Cells.Item(39, 12).Value = 1800  # L39
$ws.Cells.Item(39, 13).Value = 40.25  # M39
$ws.Cells.Item(39, 14).Value = -2392  # N39
$ws.Cells.Item(106, 8).Value = 4668.1665  # H106
$ws.Cells.Item(106, 9).Value = 2002.25  # I106
$ws.Cells.Item(106, 10).Value = 10000  # J106
$ws.Cells.Item(106, 11).Value = 2002.25  # K106
$ws.Cells.Item(106, 12).Value = 10000  # L106
$ws.Cells.Item(106, 13).Value = -1371.25  # M106
$ws.Cells.Item(106, 14).Value = -11262  # N106
$ws.Cells.Item(112, 8).Value = 550687.6  # H112
$ws.Cells.Item(112, 10).Value = 621638.5600000001  # J112
$ws.Cells.Item(112, 12).Value = 1864915.68  # L112
$ws.Cells.Item(112, 14).Value = -1867131.68  # N112
$ws.Cells.Item(138, 8).Value = 2099.34  # H138
$ws.Cells.Item(138, 9).Value = 855.5833  # I138
$ws.Cells.Item(138, 10).Value = 2492.1052  # J138
$ws.Cells.Item(138, 11).Value = 2566.7499  # K138
$ws.Cells.Item(138, 12).Value = 7476.3156  # L138
$ws.Cells.Item(138, 13).Value = 2573.2501  # M138
$ws.Cells.Item(138, 14).Value = -17756.3156  # N138

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4105.933  # H32
$ws.Cells.Item(32, 9).Value = 3171.4385  # I32
$ws.Cells.Item(32, 10).Value = 8118.7646  # J32
$ws.Cells.Item(32, 11).Value = 3171.4385  # K32
$ws.Cells.Item(32, 12).Value = 8118.7646  # L32
$ws.Cells.Item(32, 13).Value = -2884.4385  # M32
$ws.Cells.Item(32, 14).Value = -8692.7646  # N32
$ws.Cells.Item(45, 8).Value = 1168.5625  # H45
$ws.Cells.Item(45, 9).Value = 1157  # I45
$ws.Cells.Item(45, 10).Value = 1194  # J45
$ws.Cells.Item(45, 11).Value = 1157  # K45
$ws.Cells.Item(45, 12).Value = 1194  # L45
$ws.Cells.Item(45, 13).Value = -780  # M45
$ws.Cells.Item(45, 14).Value = -1948  # N45
$ws.Cells.Item(88, 8).Value = 4447084.5  # H88
$ws.Cells.Item(88, 9).Value = 5130635.5  # I88
$ws.Cells.Item(88, 10).Value = 4003.5  # J88
$ws.Cells.Item(88, 11).Value = 5130635.5  # K88
$ws.Cells.Item(88, 12).Value = 4003.5  # L88
$ws.Cells.Item(88, 13).Value = -5130229.5  # M88
$ws.Cells.Item(88, 14).Value = -4815.5  # N88
$ws.Cells.Item(91, 8).Value = 4447084.5  # H91
$ws.Cells.Item(91, 9).Value = 5130635.5  # I91
$ws.Cells.Item(91, 10).Value = 4003.5  # J91
$ws.Cells.Item(91, 11).Value = 5130635.5  # K91
$ws.Cells.Item(91, 12).Value = 4003.5  # L91
$ws.Cells.Item(91, 13).Value = -5129231.5  # M91
$ws.Cells.Item(91, 14).Value = -6811.5  # N91
$ws.Cells.Item(109, 8).Value = 30347.62  # H109
$ws.Cells.Item(109, 10).Value = 30347.62  # J109
$ws.Cells.Item(109, 12).Value = 30347.62  # L109
$ws.Cells.Item(109, 14).Value = -33121.62  # N109
$ws.Cells.Item(122, 8).Value = 3199.9092  # H122
$ws.Cells.Item(122, 9).Value = 1999.8889  # I122
$ws.Cells.Item(122, 10).Value = 8600  # J122
$ws.Cells.Item(122, 11).Value = 5999.6667  # K122
$ws.Cells.Item(122, 12).Value = 25800  # L122
$ws.Cells.Item(122, 13).Value = -3549.6667  # M122
$ws.Cells.Item(122, 14).Value = -30700  # N122

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2106.08  # H86
$ws.Cells.Item(86, 9).Value = 1885  # I86
$ws.Cells.Item(86, 10).Value = 2575.875  # J86
$ws.Cells.Item(86, 11).Value = 1885  # K86
$ws.Cells.Item(86, 12).Value = 2575.875  # L86
$ws.Cells.Item(86, 13).Value = -762  # M86
$ws.Cells.Item(86, 14).Value = -4821.875  # N86
$ws.Cells.Item(89, 8).Value = 2106.08  # H89
$ws.Cells.Item(89, 9).Value = 1885  # I89
$ws.Cells.Item(89, 10).Value = 2575.875  # J89
$ws.Cells.Item(89, 11).Value = 9425  # K89
$ws.Cells.Item(89, 12).Value = 12879.375  # L89
$ws.Cells.Item(89, 13).Value = -3809  # M89
$ws.Cells.Item(89, 14).Value = -24111.375  # N89

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 21742584  # H31
$ws.Cells.Item(31, 9).Value = 1571.0769  # I31
$ws.Cells.Item(31, 11).Value = 1571.0769  # K31
$ws.Cells.Item(31, 13).Value = -1276.0769  # M31
$ws.Cells.Item(34, 8).Value = 21742584  # H34
$ws.Cells.Item(34, 9).Value = 1571.0769  # I34
$ws.Cells.Item(34, 11).Value = 1571.0769  # K34
$ws.Cells.Item(34, 13).Value = -1369.0769  # M34
$ws.Cells.Item(58, 8).Value = 1571.8314  # H58
$ws.Cells.Item(58, 9).Value = 1354.5181  # I58
$ws.Cells.Item(58, 10).Value = 4578  # J58
$ws.Cells.Item(58, 11).Value = 1354.5181  # K58
$ws.Cells.Item(58, 12).Value = 4578  # L58
$ws.Cells.Item(58, 13).Value = -1151.5181  # M58
$ws.Cells.Item(58, 14).Value = -4984  # N58
$ws.Cells.Item(136, 8).Value = 1571.8314  # H136
$ws.Cells.Item(136, 9).Value = 1354.5181  # I136
$ws.Cells.Item(136, 10).Value = 4578  # J136
$ws.Cells.Item(136, 11).Value = 4063.5543  # K136
$ws.Cells.Item(136, 12).Value = 13734  # L136
$ws.Cells.Item(136, 13).Value = -1513.5543  # M136
$ws.Cells.Item(136, 14).Value = -18834  # N136

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(51, 8).Value = 2932.6667  # H51
$ws.Cells.Item(51, 9).Value = 899  # I51
$ws.Cells.Item(51, 11).Value = 2697  # K51
$ws.Cells.Item(51, 13).Value = -2237  # M51
$ws.Cells.Item(97, 8).Value = 0  # H97
$ws.Cells.Item(97, 9).Value = 0  # I97
$ws.Cells.Item(97, 10).Value = 0  # J97
$ws.Cells.Item(97, 11).Value = 0  # K97
$ws.Cells.Item(97, 12).ClearContents()  # L97
$ws.Cells.Item(97, 13).ClearContents()  # M97
$ws.Cells.Item(97, 14).Value = 0  # N97
$ws.Cells.Item(131, 8).Value = 7353812  # H131
$ws.Cells.Item(131, 10).Value = 933.16394  # J131
$ws.Cells.Item(131, 12).Value = 2799.49182  # L131
$ws.Cells.Item(131, 14).Value = -12879.49182  # N131

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 1792.0646  # H102
$ws.Cells.Item(102, 9).Value = 1447.5  # I102
$ws.Cells.Item(102, 10).Value = 2634.3333  # J102
$ws.Cells.Item(102, 11).Value = 1447.5  # K102
$ws.Cells.Item(102, 12).Value = 2634.3333  # L102
$ws.Cells.Item(102, 13).Value = 174.5  # M102
$ws.Cells.Item(102, 14).Value = -5878.3333  # N102

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(39, 8).Value = 0  # H39
$ws.Cells.Item(39, 10).Value = 0  # J39
$ws.Cells.Item(39, 12).Value = 0  # L39
$ws.Cells.Item(39, 14).Value = 0  # N39
$ws.Cells.Item(46, 8).Value = 3528  # H46
$ws.Cells.Item(46, 9).Value = 4000  # I46
$ws.Cells.Item(46, 10).Value = 3150.4  # J46
$ws.Cells.Item(46, 11).Value = 4000  # K46
$ws.Cells.Item(46, 12).Value = 3150.4  # L46
$ws.Cells.Item(46, 13).Value = -3812  # M46
$ws.Cells.Item(46, 14).Value = -3526.4  # N46
$ws.Cells.Item(132, 8).Value = 12557.464  # H132
$ws.Cells.Item(132, 9).Value = 11387.917  # I132
$ws.Cells.Item(132, 10).Value = 19574.75  # J132
$ws.Cells.Item(132, 11).Value = 34163.751  # K132
$ws.Cells.Item(132, 12).Value = 58724.25  # L132
$ws.Cells.Item(132, 13).Value = -31633.751  # M132
$ws.Cells.Item(132, 14).Value = -63784.25  # N132

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3074.923  # H136
$ws.Cells.Item(136, 9).Value = 2052.182  # I136
$ws.Cells.Item(136, 10).Value = 8700  # J136
$ws.Cells.Item(136, 11).Value = 6156.545999999999  # K136
$ws.Cells.Item(136, 12).Value = 26100  # L136
$ws.Cells.Item(136, 13).Value = -3606.545999999999  # M136
$ws.Cells.Item(136, 14).Value = -31200  # N136
